$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the data columns to document what each column actually holds:
#  - column B holds a numeric value with a decimal portion (21212.4444)
#  - column C holds a date value (formatted via a date number format)
# Setting A2 first (back to its own text) preserves the existing shared-string
# slot ordering, then C1/B1 are written in an order that makes the new
# shared strings land as "Date Data" (idx 2) followed by
# "Numeric Data With Decimal" (idx 3), matching the saved workbook.
$ws.Range("A2").Value = "a"
$ws.Range("C1").Value = "Date Data"
$ws.Range("B1").Value = "Numeric Data With Decimal"

# Move the active selection from C1 to B1.
$ws.Range("B1").Select()
